# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de):
#  - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - Populate "Latest Target File" (E) and "Latest Handback File" (F) columns
#    for the two real source rows (2 and 3), mirroring the Source File Name
#    (A) / Latest Handoff File (C) hyperlinked file names.
#  - Stamp "Latest Handback DateTime" (G) with the handback timestamp.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # OLE (BGR) form of RGB 6495ED, matching the workbook's existing HyperLink font color

function Set-HandbackRow($ws, $row, $mdDisplay, $mdUrl, $xlfDisplay, $xlfUrl, $handbackDateTime) {
    # Status: handed back, in sync with the English source
    $ws.Range("B$row").Value = "Handed back: in sync with en-US"

    # Latest Target File (E) - same file handed back as the target
    $ws.Hyperlinks.Add($ws.Range("E$row"), $mdUrl, "", "", $mdDisplay)
    $ws.Range("E$row").Font.Underline = 2
    $ws.Range("E$row").Font.Color = $hyperlinkColor

    # Latest Handback File (F) - the handed-back translated file
    $ws.Hyperlinks.Add($ws.Range("F$row"), $xlfUrl, "", "", $xlfDisplay)
    $ws.Range("F$row").Font.Underline = 2
    $ws.Range("F$row").Font.Color = $hyperlinkColor

    # Latest Handback DateTime (G)
    $ws.Range("G$row").Value = $handbackDateTime
}

# ---- Overview sheet mirrors the same Status text for both languages ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhMdDisplay = "617b5400-e27b-4269-bf2c-0532877aa549.md"
$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/59afd9652ef9afca00b8dbebcfc660d677d8bb54/e2e/617b5400-e27b-4269-bf2c-0532877aa549.md"
$zhXlfDisplay = "617b5400-e27b-4269-bf2c-0532877aa549.3485b3f130250f9670cc4c318f4fa47dd7fafe79.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c7a0a56d8701bf3c3e42688547dfabd78dfbdef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/617b5400-e27b-4269-bf2c-0532877aa549.3485b3f130250f9670cc4c318f4fa47dd7fafe79.zh-cn.xlf"
$zhHandback = "2016-03-04 11:07:48"

Set-HandbackRow $wsZh 2 $zhMdDisplay $zhMdUrl $zhXlfDisplay $zhXlfUrl $zhHandback
Set-HandbackRow $wsZh 3 $zhMdDisplay $zhMdUrl $zhXlfDisplay $zhXlfUrl $zhHandback

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$deMdDisplay = "617b5400-e27b-4269-bf2c-0532877aa549.md"
$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/59afd9652ef9afca00b8dbebcfc660d677d8bb54/e2e/617b5400-e27b-4269-bf2c-0532877aa549.md"
$deXlfDisplay = "617b5400-e27b-4269-bf2c-0532877aa549.3485b3f130250f9670cc4c318f4fa47dd7fafe79.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/959f3aed9ecc0b349efbffe882935b61919132e2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/617b5400-e27b-4269-bf2c-0532877aa549.3485b3f130250f9670cc4c318f4fa47dd7fafe79.de-de.xlf"
$deHandback = "2016-03-04 11:08:13"

Set-HandbackRow $wsDe 2 $deMdDisplay $deMdUrl $deXlfDisplay $deXlfUrl $deHandback
Set-HandbackRow $wsDe 3 $deMdDisplay $deMdUrl $deXlfDisplay $deXlfUrl $deHandback

Write-Output "Handback report generated for zh-cn and de-de sheets"
